$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version + Date, insert a new "Jurisdiction" row ---
$meta = $wb.Worksheets.Item("Metadata")

# Update version string (row 3, col B)
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update date string (row 8, col B)
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row right after "Contact" (row 10) / before "Description" (row 11)
# so that it becomes the new row 11: "Jurisdiction" with an empty value.
$meta.Rows.Item(11).Insert()

# Copy formatting from the row above (Contact, now row 10) so the new row matches
# the existing look (borders / wrap-text / vertical alignment) instead of getting
# a newly-synthesized style.
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""
$meta.Range("B11").ClearContents()

# --- Sheet "Elements": record the II-1 constraint on Participant1.typeId (row 5) ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Cells.Item(5, 36).Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
